$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the "base" color for the rainbow palette (row 3)
$ws.Range("B3").Value = "#1B3037"

# Correct the "complementary" column so it matches the "text" column value
# (rainbow / hotcold / coty rows had wrong complementary colors)
$ws.Range("G3").Value = $ws.Range("C3").Value()
$ws.Range("G5").Value = $ws.Range("C5").Value()
$ws.Range("G7").Value = $ws.Range("C7").Value()

# Update the active selection to C3
$ws.Range("C3").Select()
